$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,4).Value = "¥1,019,638.10"
$ws.Cells.Item(2,5).Value = "¥+19,638.10"
Set-TextValue $ws.Cells.Item(2,6) "+1.96%"
Set-TextValue $ws.Cells.Item(2,7) "+35.84%"
$ws.Cells.Item(2,8).Value = 4.648
Set-TextValue $ws.Cells.Item(2,10) "73.3%"
Set-TextValue $ws.Cells.Item(2,11) "0.1306%"
Set-TextValue $ws.Cells.Item(2,12) "6.6552%"
$ws.Cells.Item(2,13).Value = 16
$ws.Cells.Item(2,14).Value = 16
Set-TextValue $ws.Cells.Item(2,16) "20260109"
$ws.Cells.Item(3,4).Value = "¥1,042,439.98"
$ws.Cells.Item(3,5).Value = "¥+42,439.98"
Set-TextValue $ws.Cells.Item(3,6) "+4.24%"
Set-TextValue $ws.Cells.Item(3,7) "+92.44%"
$ws.Cells.Item(3,8).Value = 6.638
Set-TextValue $ws.Cells.Item(3,10) "66.7%"
Set-TextValue $ws.Cells.Item(3,11) "0.2796%"
Set-TextValue $ws.Cells.Item(3,12) "10.3159%"
$ws.Cells.Item(3,13).Value = 16
$ws.Cells.Item(3,14).Value = 16
Set-TextValue $ws.Cells.Item(3,16) "20260109"
$ws.Cells.Item(4,4).Value = "¥1,014,809.18"
$ws.Cells.Item(4,5).Value = "¥+14,809.18"
Set-TextValue $ws.Cells.Item(4,6) "+1.48%"
Set-TextValue $ws.Cells.Item(4,7) "+26.05%"
$ws.Cells.Item(4,8).Value = 4.431
Set-TextValue $ws.Cells.Item(4,10) "53.3%"
Set-TextValue $ws.Cells.Item(4,11) "0.0986%"
Set-TextValue $ws.Cells.Item(4,12) "5.1597%"
$ws.Cells.Item(4,13).Value = 16
$ws.Cells.Item(4,14).Value = 16
Set-TextValue $ws.Cells.Item(4,16) "20260109"
$ws.Cells.Item(5,4).Value = "¥1,003,464.49"
$ws.Cells.Item(5,5).Value = "¥+3,464.49"
Set-TextValue $ws.Cells.Item(5,6) "+0.35%"
Set-TextValue $ws.Cells.Item(5,7) "+5.60%"
$ws.Cells.Item(5,8).Value = 1.575
Set-TextValue $ws.Cells.Item(5,10) "53.3%"
Set-TextValue $ws.Cells.Item(5,11) "0.0232%"
Set-TextValue $ws.Cells.Item(5,12) "2.4506%"
$ws.Cells.Item(5,13).Value = 16
$ws.Cells.Item(5,14).Value = 16
Set-TextValue $ws.Cells.Item(5,16) "20260109"
$ws.Cells.Item(6,4).Value = "¥1,001,518.00"
$ws.Cells.Item(6,5).Value = "¥+1,518.00"
Set-TextValue $ws.Cells.Item(6,6) "+0.15%"
Set-TextValue $ws.Cells.Item(6,7) "+2.42%"
$ws.Cells.Item(6,8).Value = 0.472
Set-TextValue $ws.Cells.Item(6,10) "53.3%"
Set-TextValue $ws.Cells.Item(6,11) "0.0101%"
Set-TextValue $ws.Cells.Item(6,12) "1.2193%"
$ws.Cells.Item(6,13).Value = 16
$ws.Cells.Item(6,14).Value = 16
Set-TextValue $ws.Cells.Item(6,16) "20260109"
$ws.Cells.Item(7,4).Value = "¥1,040,563.19"
$ws.Cells.Item(7,5).Value = "¥+40,563.19"
Set-TextValue $ws.Cells.Item(7,6) "+4.06%"
Set-TextValue $ws.Cells.Item(7,7) "+87.06%"
$ws.Cells.Item(7,8).Value = 8.3
Set-TextValue $ws.Cells.Item(7,10) "66.7%"
Set-TextValue $ws.Cells.Item(7,11) "0.2667%"
Set-TextValue $ws.Cells.Item(7,12) "7.8577%"
$ws.Cells.Item(7,13).Value = 16
$ws.Cells.Item(7,14).Value = 16
Set-TextValue $ws.Cells.Item(7,16) "20260109"
$ws.Cells.Item(8,4).Value = "¥1,026,881.61"
$ws.Cells.Item(8,5).Value = "¥+26,881.61"
Set-TextValue $ws.Cells.Item(8,6) "+2.69%"
Set-TextValue $ws.Cells.Item(8,7) "+51.86%"
$ws.Cells.Item(8,8).Value = 6.192
Set-TextValue $ws.Cells.Item(8,10) "53.3%"
Set-TextValue $ws.Cells.Item(8,11) "0.1779%"
Set-TextValue $ws.Cells.Item(8,12) "6.9223%"
$ws.Cells.Item(8,13).Value = 16
$ws.Cells.Item(8,14).Value = 16
Set-TextValue $ws.Cells.Item(8,16) "20260109"
$ws.Cells.Item(9,4).Value = "¥1,030,262.17"
$ws.Cells.Item(9,5).Value = "¥+30,262.17"
Set-TextValue $ws.Cells.Item(9,6) "+3.03%"
Set-TextValue $ws.Cells.Item(9,7) "+59.93%"
$ws.Cells.Item(9,8).Value = 6.98
Set-TextValue $ws.Cells.Item(9,10) "60.0%"
Set-TextValue $ws.Cells.Item(9,11) "0.1999%"
Set-TextValue $ws.Cells.Item(9,12) "6.9331%"
$ws.Cells.Item(9,13).Value = 16
$ws.Cells.Item(9,14).Value = 16
Set-TextValue $ws.Cells.Item(9,16) "20260109"
$ws.Cells.Item(10,4).Value = "¥1,052,440.53"
$ws.Cells.Item(10,5).Value = "¥+52,440.53"
Set-TextValue $ws.Cells.Item(10,6) "+5.24%"
Set-TextValue $ws.Cells.Item(10,7) "+123.67%"
$ws.Cells.Item(10,8).Value = 7.436
Set-TextValue $ws.Cells.Item(10,10) "66.7%"
Set-TextValue $ws.Cells.Item(10,11) "0.3439%"
Set-TextValue $ws.Cells.Item(10,12) "11.3886%"
$ws.Cells.Item(10,13).Value = 16
$ws.Cells.Item(10,14).Value = 16
Set-TextValue $ws.Cells.Item(10,16) "20260109"
$ws.Cells.Item(11,4).Value = "¥1,002,594.00"
$ws.Cells.Item(11,5).Value = "¥+2,594.00"
Set-TextValue $ws.Cells.Item(11,6) "+0.26%"
Set-TextValue $ws.Cells.Item(11,7) "+4.16%"
$ws.Cells.Item(11,8).Value = 2.387
Set-TextValue $ws.Cells.Item(11,10) "53.3%"
Set-TextValue $ws.Cells.Item(11,11) "0.0173%"
Set-TextValue $ws.Cells.Item(11,12) "0.9959%"
$ws.Cells.Item(11,13).Value = 16
$ws.Cells.Item(11,14).Value = 16
Set-TextValue $ws.Cells.Item(11,16) "20260109"
$ws.Cells.Item(12,4).Value = "¥1,011,144.42"
$ws.Cells.Item(12,5).Value = "¥+11,144.42"
Set-TextValue $ws.Cells.Item(12,6) "+1.11%"
Set-TextValue $ws.Cells.Item(12,7) "+19.07%"
$ws.Cells.Item(12,8).Value = 3.44
Set-TextValue $ws.Cells.Item(12,10) "53.3%"
Set-TextValue $ws.Cells.Item(12,11) "0.0744%"
Set-TextValue $ws.Cells.Item(12,12) "4.8735%"
$ws.Cells.Item(12,13).Value = 16
$ws.Cells.Item(12,14).Value = 16
Set-TextValue $ws.Cells.Item(12,16) "20260109"
$ws.Cells.Item(13,4).Value = "¥1,042,601.48"
$ws.Cells.Item(13,5).Value = "¥+42,601.48"
Set-TextValue $ws.Cells.Item(13,6) "+4.26%"
Set-TextValue $ws.Cells.Item(13,7) "+92.91%"
$ws.Cells.Item(13,8).Value = 6.292
Set-TextValue $ws.Cells.Item(13,10) "66.7%"
Set-TextValue $ws.Cells.Item(13,11) "0.2809%"
Set-TextValue $ws.Cells.Item(13,12) "10.9356%"
$ws.Cells.Item(13,13).Value = 16
$ws.Cells.Item(13,14).Value = 16
Set-TextValue $ws.Cells.Item(13,16) "20260109"
$ws.Cells.Item(14,4).Value = "¥1,007,811.62"
$ws.Cells.Item(14,5).Value = "¥+7,811.62"
Set-TextValue $ws.Cells.Item(14,6) "+0.78%"
Set-TextValue $ws.Cells.Item(14,7) "+13.04%"
$ws.Cells.Item(14,8).Value = 3.849
Set-TextValue $ws.Cells.Item(14,10) "53.3%"
Set-TextValue $ws.Cells.Item(14,11) "0.0521%"
Set-TextValue $ws.Cells.Item(14,12) "2.8938%"
$ws.Cells.Item(14,13).Value = 16
$ws.Cells.Item(14,14).Value = 16
Set-TextValue $ws.Cells.Item(14,16) "20260109"
$ws.Cells.Item(15,4).Value = "¥1,000,596.60"
$ws.Cells.Item(15,5).Value = "¥+596.60"
Set-TextValue $ws.Cells.Item(15,6) "+0.06%"
Set-TextValue $ws.Cells.Item(15,7) "+0.94%"
$ws.Cells.Item(15,8).Value = -2.063
Set-TextValue $ws.Cells.Item(15,10) "53.3%"
Set-TextValue $ws.Cells.Item(15,11) "0.0040%"
Set-TextValue $ws.Cells.Item(15,12) "0.4737%"
$ws.Cells.Item(15,13).Value = 16
$ws.Cells.Item(15,14).Value = 16
Set-TextValue $ws.Cells.Item(15,16) "20260109"
$ws.Cells.Item(16,4).Value = "¥1,000,663.50"
$ws.Cells.Item(16,5).Value = "¥+663.50"
Set-TextValue $ws.Cells.Item(16,6) "+0.07%"
Set-TextValue $ws.Cells.Item(16,7) "+1.05%"
$ws.Cells.Item(16,8).Value = -1.82
Set-TextValue $ws.Cells.Item(16,10) "53.3%"
Set-TextValue $ws.Cells.Item(16,11) "0.0044%"
Set-TextValue $ws.Cells.Item(16,12) "0.4752%"
$ws.Cells.Item(16,13).Value = 16
$ws.Cells.Item(16,14).Value = 16
Set-TextValue $ws.Cells.Item(16,16) "20260109"

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,4).Value = "¥1,019,638.10"
$ws.Cells.Item(2,5).Value = "¥+19,638.10"
Set-TextValue $ws.Cells.Item(2,6) "+1.96%"
Set-TextValue $ws.Cells.Item(2,7) "+35.84%"
$ws.Cells.Item(2,8).Value = 4.648
Set-TextValue $ws.Cells.Item(2,10) "73.3%"
Set-TextValue $ws.Cells.Item(2,11) "0.1306%"
Set-TextValue $ws.Cells.Item(2,12) "6.6552%"
$ws.Cells.Item(2,13).Value = 16
$ws.Cells.Item(2,14).Value = 16
Set-TextValue $ws.Cells.Item(2,16) "20260109"
$ws.Cells.Item(3,4).Value = "¥1,042,439.98"
$ws.Cells.Item(3,5).Value = "¥+42,439.98"
Set-TextValue $ws.Cells.Item(3,6) "+4.24%"
Set-TextValue $ws.Cells.Item(3,7) "+92.44%"
$ws.Cells.Item(3,8).Value = 6.638
Set-TextValue $ws.Cells.Item(3,10) "66.7%"
Set-TextValue $ws.Cells.Item(3,11) "0.2796%"
Set-TextValue $ws.Cells.Item(3,12) "10.3159%"
$ws.Cells.Item(3,13).Value = 16
$ws.Cells.Item(3,14).Value = 16
Set-TextValue $ws.Cells.Item(3,16) "20260109"
$ws.Cells.Item(4,4).Value = "¥1,014,809.18"
$ws.Cells.Item(4,5).Value = "¥+14,809.18"
Set-TextValue $ws.Cells.Item(4,6) "+1.48%"
Set-TextValue $ws.Cells.Item(4,7) "+26.05%"
$ws.Cells.Item(4,8).Value = 4.431
Set-TextValue $ws.Cells.Item(4,10) "53.3%"
Set-TextValue $ws.Cells.Item(4,11) "0.0986%"
Set-TextValue $ws.Cells.Item(4,12) "5.1597%"
$ws.Cells.Item(4,13).Value = 16
$ws.Cells.Item(4,14).Value = 16
Set-TextValue $ws.Cells.Item(4,16) "20260109"
$ws.Cells.Item(5,4).Value = "¥1,003,464.49"
$ws.Cells.Item(5,5).Value = "¥+3,464.49"
Set-TextValue $ws.Cells.Item(5,6) "+0.35%"
Set-TextValue $ws.Cells.Item(5,7) "+5.60%"
$ws.Cells.Item(5,8).Value = 1.575
Set-TextValue $ws.Cells.Item(5,10) "53.3%"
Set-TextValue $ws.Cells.Item(5,11) "0.0232%"
Set-TextValue $ws.Cells.Item(5,12) "2.4506%"
$ws.Cells.Item(5,13).Value = 16
$ws.Cells.Item(5,14).Value = 16
Set-TextValue $ws.Cells.Item(5,16) "20260109"
$ws.Cells.Item(6,4).Value = "¥1,001,518.00"
$ws.Cells.Item(6,5).Value = "¥+1,518.00"
Set-TextValue $ws.Cells.Item(6,6) "+0.15%"
Set-TextValue $ws.Cells.Item(6,7) "+2.42%"
$ws.Cells.Item(6,8).Value = 0.472
Set-TextValue $ws.Cells.Item(6,10) "53.3%"
Set-TextValue $ws.Cells.Item(6,11) "0.0101%"
Set-TextValue $ws.Cells.Item(6,12) "1.2193%"
$ws.Cells.Item(6,13).Value = 16
$ws.Cells.Item(6,14).Value = 16
Set-TextValue $ws.Cells.Item(6,16) "20260109"

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,4).Value = "¥1,040,563.19"
$ws.Cells.Item(2,5).Value = "¥+40,563.19"
Set-TextValue $ws.Cells.Item(2,6) "+4.06%"
Set-TextValue $ws.Cells.Item(2,7) "+87.06%"
$ws.Cells.Item(2,8).Value = 8.3
Set-TextValue $ws.Cells.Item(2,10) "66.7%"
Set-TextValue $ws.Cells.Item(2,11) "0.2667%"
Set-TextValue $ws.Cells.Item(2,12) "7.8577%"
$ws.Cells.Item(2,13).Value = 16
$ws.Cells.Item(2,14).Value = 16
Set-TextValue $ws.Cells.Item(2,16) "20260109"
$ws.Cells.Item(3,4).Value = "¥1,026,881.61"
$ws.Cells.Item(3,5).Value = "¥+26,881.61"
Set-TextValue $ws.Cells.Item(3,6) "+2.69%"
Set-TextValue $ws.Cells.Item(3,7) "+51.86%"
$ws.Cells.Item(3,8).Value = 6.192
Set-TextValue $ws.Cells.Item(3,10) "53.3%"
Set-TextValue $ws.Cells.Item(3,11) "0.1779%"
Set-TextValue $ws.Cells.Item(3,12) "6.9223%"
$ws.Cells.Item(3,13).Value = 16
$ws.Cells.Item(3,14).Value = 16
Set-TextValue $ws.Cells.Item(3,16) "20260109"
$ws.Cells.Item(4,4).Value = "¥1,030,262.17"
$ws.Cells.Item(4,5).Value = "¥+30,262.17"
Set-TextValue $ws.Cells.Item(4,6) "+3.03%"
Set-TextValue $ws.Cells.Item(4,7) "+59.93%"
$ws.Cells.Item(4,8).Value = 6.98
Set-TextValue $ws.Cells.Item(4,10) "60.0%"
Set-TextValue $ws.Cells.Item(4,11) "0.1999%"
Set-TextValue $ws.Cells.Item(4,12) "6.9331%"
$ws.Cells.Item(4,13).Value = 16
$ws.Cells.Item(4,14).Value = 16
Set-TextValue $ws.Cells.Item(4,16) "20260109"
$ws.Cells.Item(5,4).Value = "¥1,052,440.53"
$ws.Cells.Item(5,5).Value = "¥+52,440.53"
Set-TextValue $ws.Cells.Item(5,6) "+5.24%"
Set-TextValue $ws.Cells.Item(5,7) "+123.67%"
$ws.Cells.Item(5,8).Value = 7.436
Set-TextValue $ws.Cells.Item(5,10) "66.7%"
Set-TextValue $ws.Cells.Item(5,11) "0.3439%"
Set-TextValue $ws.Cells.Item(5,12) "11.3886%"
$ws.Cells.Item(5,13).Value = 16
$ws.Cells.Item(5,14).Value = 16
Set-TextValue $ws.Cells.Item(5,16) "20260109"
$ws.Cells.Item(6,4).Value = "¥1,002,594.00"
$ws.Cells.Item(6,5).Value = "¥+2,594.00"
Set-TextValue $ws.Cells.Item(6,6) "+0.26%"
Set-TextValue $ws.Cells.Item(6,7) "+4.16%"
$ws.Cells.Item(6,8).Value = 2.387
Set-TextValue $ws.Cells.Item(6,10) "53.3%"
Set-TextValue $ws.Cells.Item(6,11) "0.0173%"
Set-TextValue $ws.Cells.Item(6,12) "0.9959%"
$ws.Cells.Item(6,13).Value = 16
$ws.Cells.Item(6,14).Value = 16
Set-TextValue $ws.Cells.Item(6,16) "20260109"

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2,4).Value = "¥1,011,144.42"
$ws.Cells.Item(2,5).Value = "¥+11,144.42"
Set-TextValue $ws.Cells.Item(2,6) "+1.11%"
Set-TextValue $ws.Cells.Item(2,7) "+19.07%"
$ws.Cells.Item(2,8).Value = 3.44
Set-TextValue $ws.Cells.Item(2,10) "53.3%"
Set-TextValue $ws.Cells.Item(2,11) "0.0744%"
Set-TextValue $ws.Cells.Item(2,12) "4.8735%"
$ws.Cells.Item(2,13).Value = 16
$ws.Cells.Item(2,14).Value = 16
Set-TextValue $ws.Cells.Item(2,16) "20260109"
$ws.Cells.Item(3,4).Value = "¥1,042,601.48"
$ws.Cells.Item(3,5).Value = "¥+42,601.48"
Set-TextValue $ws.Cells.Item(3,6) "+4.26%"
Set-TextValue $ws.Cells.Item(3,7) "+92.91%"
$ws.Cells.Item(3,8).Value = 6.292
Set-TextValue $ws.Cells.Item(3,10) "66.7%"
Set-TextValue $ws.Cells.Item(3,11) "0.2809%"
Set-TextValue $ws.Cells.Item(3,12) "10.9356%"
$ws.Cells.Item(3,13).Value = 16
$ws.Cells.Item(3,14).Value = 16
Set-TextValue $ws.Cells.Item(3,16) "20260109"
$ws.Cells.Item(4,4).Value = "¥1,007,811.62"
$ws.Cells.Item(4,5).Value = "¥+7,811.62"
Set-TextValue $ws.Cells.Item(4,6) "+0.78%"
Set-TextValue $ws.Cells.Item(4,7) "+13.04%"
$ws.Cells.Item(4,8).Value = 3.849
Set-TextValue $ws.Cells.Item(4,10) "53.3%"
Set-TextValue $ws.Cells.Item(4,11) "0.0521%"
Set-TextValue $ws.Cells.Item(4,12) "2.8938%"
$ws.Cells.Item(4,13).Value = 16
$ws.Cells.Item(4,14).Value = 16
Set-TextValue $ws.Cells.Item(4,16) "20260109"
$ws.Cells.Item(5,4).Value = "¥1,000,596.60"
$ws.Cells.Item(5,5).Value = "¥+596.60"
Set-TextValue $ws.Cells.Item(5,6) "+0.06%"
Set-TextValue $ws.Cells.Item(5,7) "+0.94%"
$ws.Cells.Item(5,8).Value = -2.063
Set-TextValue $ws.Cells.Item(5,10) "53.3%"
Set-TextValue $ws.Cells.Item(5,11) "0.0040%"
Set-TextValue $ws.Cells.Item(5,12) "0.4737%"
$ws.Cells.Item(5,13).Value = 16
$ws.Cells.Item(5,14).Value = 16
Set-TextValue $ws.Cells.Item(5,16) "20260109"
$ws.Cells.Item(6,4).Value = "¥1,000,663.50"
$ws.Cells.Item(6,5).Value = "¥+663.50"
Set-TextValue $ws.Cells.Item(6,6) "+0.07%"
Set-TextValue $ws.Cells.Item(6,7) "+1.05%"
$ws.Cells.Item(6,8).Value = -1.82
Set-TextValue $ws.Cells.Item(6,10) "53.3%"
Set-TextValue $ws.Cells.Item(6,11) "0.0044%"
Set-TextValue $ws.Cells.Item(6,12) "0.4752%"
$ws.Cells.Item(6,13).Value = 16
$ws.Cells.Item(6,14).Value = 16
Set-TextValue $ws.Cells.Item(6,16) "20260109"
